$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.504.47"
$ws.Range("E2").Value = "  -2.38%  "

$ws.Range("D3").Value = "1.997.04"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.34"
$ws.Range("E5").Value = "  -9.53%  "

$ws.Range("E6").Value = "  -2.83%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.79"
$ws.Range("E8").Value = "  -2.32%  "

$ws.Range("E9").Value = "  -3.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.42"
$ws.Range("E10").Value = "  +3.25%  "

$ws.Range("E11").Value = "  -2.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0983"
$ws.Range("E12").Value = "  -3.30%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.290.93"
$ws.Range("E13").Value = "  -0.52%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.18"
$ws.Range("E14").Value = "  -0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.29"
$ws.Range("E15").Value = "  -3.25%  "

$ws.Range("E16").Value = "  -5.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.07"
$ws.Range("E17").Value = "  -3.12%  "

$ws.Range("D18").Value = "2.005.58"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("D19").Value = "36.583.09"
$ws.Range("E19").Value = "  -1.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.69"
$ws.Range("E20").Value = "  -2.89%  "

$ws.Range("D21").Value = "0.0₃0803"
$ws.Range("E21").Value = "  -3.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.26"
$ws.Range("E22").Value = "  +2.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.75"
$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("E26").Value = "  -8.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.84"
$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.67"
$ws.Range("E28").Value = "  -2.10%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.38"
$ws.Range("E29").Value = "  +3.58%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  -3.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "18.77"
$ws.Range("E31").Value = "  -4.38%  "

$ws.Range("E32").Value = "  -2.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.37"
$ws.Range("E33").Value = "  -5.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("E34").Value = "  -6.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.25"
$ws.Range("E35").Value = "  -6.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  -0.83%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("E39").Value = "  -3.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.59"
$ws.Range("E40").Value = "  +4.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.99"
$ws.Range("E41").Value = "  -1.71%  "

$ws.Range("D42").Value = "1.451.14"
$ws.Range("E42").Value = "  +2.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0926"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("E44").Value = "  -4.51%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").Value = "  -8.07%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.19"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.14"
$ws.Range("E47").Value = "  -3.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.992"
$ws.Range("E48").Value = "  -3.07%  "

$ws.Range("E49").Value = "  -0.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.83"
$ws.Range("E50").Value = "  -2.74%  "

$ws.Range("E51").Value = "  +5.53%  "
